# Auto-generated edit script: update cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "34.787.00"
$c.Style = "Normal"

$ws.Range("E2").Value = "  -2.51%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.805.71"
$c.Style = "Normal"

$ws.Range("E3").Value = "  -3.08%  "

$ws.Range("E4").Value = "  +0.28%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "230.37"
$c.Style = "Normal"

$ws.Range("E5").Value = "  -0.36%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.602"
$c.Style = "Normal"

$ws.Range("E6").Value = "  -1.72%  "

$ws.Range("E7").Value = "  +0.28%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "38.62"
$c.Style = "Normal"

$ws.Range("E8").Value = "  -9.13%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.316"
$c.Style = "Normal"

$ws.Range("E9").Value = "  +2.13%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0676"
$c.Style = "Normal"

$ws.Range("E10").Value = "  -2.89%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0991"
$c.Style = "Normal"

$ws.Range("E11").Value = "  -2.13%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "2.069.52"
$c.Style = "Normal"

$ws.Range("E12").Value = "  -3.00%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.823.52"
$c.Style = "Normal"

$ws.Range("E13").Value = "  -2.18%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.656"
$c.Style = "Normal"

$ws.Range("E14").Value = "  -3.44%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "10.80"
$c.Style = "Normal"

$ws.Range("E15").Value = "  -7.13%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "4.52"
$c.Style = "Normal"

$ws.Range("E16").Value = "  -4.40%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "34.828.95"
$c.Style = "Normal"

$ws.Range("E17").Value = "  -2.45%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "69.00"
$c.Style = "Normal"

$ws.Range("E18").Value = "  -2.06%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.0₃0777"
$c.Style = "Normal"

$ws.Range("E19").Value = "  -3.42%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "238.53"
$c.Style = "Normal"

$ws.Range("E20").Value = "  -4.19%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "11.67"
$c.Style = "Normal"

$ws.Range("E21").Value = "  -5.08%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.59"
$c.Style = "Normal"

$ws.Range("E22").Value = "  -3.24%  "

$ws.Range("E23").Value = "  +0.26%  "

$ws.Range("E24").Value = "  -0.40%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "173.54"
$c.Style = "Normal"

$ws.Range("E25").Value = "  +1.71%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.70"
$c.Style = "Normal"

$ws.Range("E26").Value = "  -3.76%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "17.07"
$c.Style = "Normal"

$ws.Range("E27").Value = "  -4.75%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.118"
$c.Style = "Normal"

$ws.Range("E28").Value = "  -3.50%  "

$ws.Range("E29").Value = "  +4.92%  "

$ws.Range("E30").Value = "  +0.25%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.95"
$c.Style = "Normal"

$ws.Range("E31").Value = "  -0.26%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.0542"
$c.Style = "Normal"

$ws.Range("E32").Value = "  -0.88%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.90"
$c.Style = "Normal"

$ws.Range("E33").Value = "  -4.74%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.18"
$c.Style = "Normal"

$ws.Range("E34").Value = "  +7.75%  "

$ws.Range("E35").Value = "  -8.09%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.681"
$c.Style = "Normal"

$ws.Range("E36").Value = "  -1.49%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "90.40"
$c.Style = "Normal"

$ws.Range("E37").Value = "  -11.01%  "

$ws.Range("E38").Value = "  +6.20%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.306.12"
$c.Style = "Normal"

$ws.Range("E39").Value = "  -4.69%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0190"
$c.Style = "Normal"

$ws.Range("E40").Value = "  -3.28%  "

$ws.Range("E41").Value = "  -0.79%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.949"
$c.Style = "Normal"

$ws.Range("E42").Value = "  -6.21%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "14.02"
$c.Style = "Normal"

$ws.Range("E43").Value = "  -5.70%  "

$ws.Range("E44").Value = "  -12.53%  "

$ws.Range("E45").Value = "  -4.91%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0511"
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "6.09"
$c.Style = "Normal"

$ws.Range("E47").Value = "  -3.32%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.997.72"
$c.Style = "Normal"

$ws.Range("E48").Value = "  -1.70%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0680"
$c.Style = "Normal"

$ws.Range("E49").Value = "  +8.27%  "

$ws.Range("E50").Value = "  +0.23%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "98.31"
$c.Style = "Normal"

$ws.Range("E51").Value = "  -6.06%  "
